$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.553094
$ws.Range("H2").Value = 4.659282
$ws.Range("I2").Value = 0.6859765954652609
$ws.Range("J2").Value = 0.6859765954652609
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("Q2").Value = 19.190142322164
$ws.Range("R2").Value = 172.711280899476
$ws.Range("S2").Value = 0.6859765954652609
$ws.Range("T2").Value = 0.6859765954652609

# Row 3 updates
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7109686666666667
$ws.Range("H3").Value = 2.132906
$ws.Range("I3").Value = 0.314023404534739
$ws.Range("J3").Value = 0.314023404534739
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.35607266666667
$ws.Range("N3").Value = 37.068218
$ws.Range("Q3").Value = 8.784780509056446
$ws.Range("R3").Value = 79.06302458150802
$ws.Range("S3").Value = 0.314023404534739
$ws.Range("T3").Value = 0.314023404534739
